$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the existing student marks data so only the new header row remains
$ws.UsedRange.Clear()

# Write the new header row with all student detail columns
$ws.Range("A1").Value = "Reg No"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "GitHub Username"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Contact No"
$ws.Range("F1").Value = "Birthday"
$ws.Range("G1").Value = "Gender"
